# Apply the commit "modify readme and add rhat files" to the Author
# Contributions Checklist form:
#   1. Check 7 specific checkboxes (w14:checked 0 -> 1, glyph ☐ -> ☒)
#   2. Add a sentence about running scripts into the first (empty) bordered
#      paragraph under the "Instructions" heading.

$d = $word.ActiveDocument

# --- 1. Tick the seven checkbox content controls -------------------------
$checkedGlyph = [char]0x2612   # ☒ BALLOT BOX WITH X

$targetIds = @(
    "-1624604447",   # Scope: Any numbers provided in text in the paper
    "1516654062",    # Scope: The computational method(s) presented in the paper ...
    "-1947686302",   # Scope: All tables and figures in the paper
    "-1264292286",   # Workflow details: Single master code file
    "-1899509192",   # Workflow details: Wrapper (shell) script(s)
    "-1145584448",   # Workflow details: Other (more detail in 'Instructions' below)
    "1316299814"      # Expected run-time: >8 hours
)

foreach ($cc in $d.ContentControls) {
    if ($cc.Type -eq 8) {
        $idstr = [string]$cc.ID
        if ($targetIds -contains $idstr) {
            $cc.Checked = $true
            $cc.Range.Text = $checkedGlyph
        }
    }
}

# --- 2. Add the "Every script in code/..." sentence under Instructions ---
# The paragraph immediately following the "Instructions" heading (inside the
# bordered placeholder sdt) is empty; locate it by scanning the document's
# paragraphs for the one right after the "Instructions" heading paragraph.

$paraCount = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($txt -eq "Instructions") {
        $targetIndex = $i + 1
        break
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)
    $start = $target.Range.Start

    $text1 = "Every script in code/framingham/, code/leukemia/, code/veteran and code/sims can be executed by clicking on " + [char]0x201C + "r"
    $text2 = "un" + [char]0x201D
    $text3 = " and reproduce plots and numbers mentioned in the paper."

    # Each InsertAfter on a collapsed range inserts immediately at that fixed
    # point (ahead of whatever follows), so insert the three chunks in
    # reverse order at the same anchor to end up with the correct final
    # left-to-right reading order: text1 + text2 + text3.
    $r3 = $d.Range($start, $start)
    $r3.InsertAfter($text3)

    $r2 = $d.Range($start, $start)
    $r2.InsertAfter($text2)

    $r1 = $d.Range($start, $start)
    $r1.InsertAfter($text1)
}
